$wb = $excel.ActiveWorkbook

# Rename the BOM worksheet to "Stock" (pricing doc revamped w/ variants).
$ws = $wb.Worksheets.Item("ovrBeaconGateway-wifi_BOM")
$ws.Name = "Stock"
